$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text in the source sheet (e.g. "7.90", "3.295.18").
# Assigning a numeric-looking string via .Value to a General-formatted cell would make
# Excel coerce it to a real number (and round-trip floats like 556.82 -> 556.82000000000005),
# so each Price cell is briefly marked as Text, written, then restored to the Normal style
# (removing the temporary format) to keep the cell visually/structurally like the original.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '60.256.22'
$ws.Range("E2").Value = '  -2.83%  '
Set-TextValue "D3" '3.298.02'
$ws.Range("E3").Value = '  -3.59%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '556.82'
Set-TextValue "D6" '140.81'
$ws.Range("E6").Value = '  -8.57%  '
$ws.Range("E7").Value = '  -0.02%  '
Set-TextValue "D8" '3.296.34'
$ws.Range("E8").Value = '  -3.67%  '
$ws.Range("E9").Value = '  -3.63%  '
Set-TextValue "D10" '7.92'
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("E11").Value = '  -5.18%  '
$ws.Range("E12").Value = '  -2.45%  '
Set-TextValue "D13" '3.860.28'
$ws.Range("E13").Value = '  -3.60%  '
Set-TextValue "D15" '26.61'
$ws.Range("E15").Value = '  -6.16%  '
Set-TextValue "D16" '3.295.34'
$ws.Range("E16").Value = '  -3.56%  '
$ws.Range("E17").Value = '  -5.15%  '
Set-TextValue "D18" '60.240.07'
$ws.Range("E18").Value = '  -2.91%  '
Set-TextValue "D19" '6.05'
$ws.Range("E19").Value = '  -7.87%  '
$ws.Range("E20").Value = '  -5.37%  '
Set-TextValue "D21" '8.49'
$ws.Range("E21").Value = '  -5.19%  '
Set-TextValue "D22" '373.43'
$ws.Range("E22").Value = '  -2.28%  '
Set-TextValue "D23" '72.67'
$ws.Range("E23").Value = '  -4.30%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  -6.98%  '
Set-TextValue "D26" '3.429.77'
$ws.Range("E26").Value = '  -3.71%  '
$ws.Range("E27").Value = '  -9.82%  '
$ws.Range("E28").Value = '  -2.61%  '
Set-TextValue "D29" '0.999'
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  -8.25%  '
$ws.Range("E31").Value = '  -0.02%  '
Set-TextValue "D32" '2.02'
$ws.Range("E32").Value = '  -5.28%  '
Set-TextValue "D33" '7.40'
$ws.Range("E33").Value = '  -6.07%  '
Set-TextValue "D34" '22.51'
$ws.Range("E34").Value = '  -3.39%  '
$ws.Range("E35").Value = '  -7.86%  '
Set-TextValue "D36" '165.97'
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("E37").Value = '  -9.68%  '
$ws.Range("E38").Value = '  -5.12%  '
$ws.Range("E39").Value = '  -5.01%  '
Set-TextValue "D40" '3.325.95'
$ws.Range("E40").Value = '  -3.80%  '
$ws.Range("E41").Value = '  -8.26%  '
Set-TextValue "D42" '25.32'
$ws.Range("E42").Value = '  -18.40%  '
Set-TextValue "D43" '41.62'
$ws.Range("E44").Value = '  -4.28%  '
$ws.Range("E45").Value = '  -4.53%  '
$ws.Range("E46").Value = '  -7.93%  '
$ws.Range("E47").Value = '  -6.86%  '
$ws.Range("E48").Value = '  -0.05%  '
Set-TextValue "D49" '2.318.50'
$ws.Range("E49").Value = '  -9.39%  '
Set-TextValue "D50" '21.46'
$ws.Range("E50").Value = '  -7.05%  '
$ws.Range("E51").Value = '  -7.24%  '
